{"js": "// Update the date title and the 25 populated answer cells in the\n// two-digit-divided-by-one-digit-number worksheet table.\n//\n// The table has 20 rows x 5 columns; only every 4th row (0, 4, 8, 12, 16)\n// actually has text (the rows between are intentionally blank spacer\n// rows). Replacement is strictly positional (row, col) -> new text,\n// because several original cell values repeat verbatim elsewhere in the\n// table with *different* replacements, so a global find/replace would be\n// ambiguous.\n\nconst body = context.document.body;\n\n// ---- 1. Title paragraph -------------------------------------------------\n// Paragraph.text is read-only in the Word JS API, so replace via the\n// paragraph's own range instead of assigning the property directly.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  paragraphs.items[0].getRange().insertText(\"2025-10-29 Wednesday\", \"Replace\");\n}\n\n// ---- 2. Table cells ------------------------------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (rowIndex, colIndex) -> new text, in document order.\nconst replacements = [\n  [0, 0, \"10\u00f77=1, 3\"],\n  [0, 1, \"47\u00f79=5, 2\"],\n  [0, 2, \"54\u00f76=9, 0\"],\n  [0, 3, \"24\u00f74=6, 0\"],\n  [0, 4, \"64\u00f78=8, 0\"],\n\n  [4, 0, \"71\u00f79=7, 8\"],\n  [4, 1, \"32\u00f77=4, 4\"],\n  [4, 2, \"85\u00f75=17, 0\"],\n  [4, 3, \"86\u00f78=10, 6\"],\n  [4, 4, \"37\u00f76=6, 1\"],\n\n  [8, 0, \"44\u00f72=22, 0\"],\n  [8, 1, \"34\u00f75=6, 4\"],\n  [8, 2, \"27\u00f72=13, 1\"],\n  [8, 3, \"17\u00f77=2, 3\"],\n  [8, 4, \"34\u00f72=17, 0\"],\n\n  [12, 0, \"93\u00f77=13, 2\"],\n  [12, 1, \"86\u00f74=21, 2\"],\n  [12, 2, \"82\u00f78=10, 2\"],\n  [12, 3, \"13\u00f79=1, 4\"],\n  [12, 4, \"52\u00f74=13, 0\"],\n\n  [16, 0, \"86\u00f74=21, 2\"],\n  [16, 1, \"47\u00f72=23, 1\"],\n  [16, 2, \"70\u00f75=14, 0\"],\n  [16, 3, \"16\u00f78=2, 0\"],\n  [16, 4, \"66\u00f77=9, 3\"],\n];\n\nfor (const [row, col, text] of replacements) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and the 25 populated answer cells in the\n# two-digit-divided-by-one-digit-number worksheet table.\n#\n# The table has 20 rows x 5 columns; only every 4th row (COM rows 1, 5, 9,\n# 13, 17 - 1-based) actually has text, the rows between are intentionally\n# blank spacer rows. Replacement is strictly positional (row, col) -> new\n# text, because several original cell values repeat verbatim elsewhere in\n# the table with *different* replacements, so a global find/replace would\n# be ambiguous.\n\n$d = $word.ActiveDocument\n\n# ---- 1. Title paragraph --------------------------------------------------\n$d.Paragraphs.Item(1).Range.Text = \"2025-10-29 Wednesday\"\n\n# ---- 2. Table cells -------------------------------------------------------\n$t = $d.Tables.Item(1)\n\n# row, col (both 1-based, matching Word COM's Table.Cell(row, col)) -> text\n$replacements = @(\n    @(1, 1, \"10\u00f77=1, 3\"),\n    @(1, 2, \"47\u00f79=5, 2\"),\n    @(1, 3, \"54\u00f76=9, 0\"),\n    @(1, 4, \"24\u00f74=6, 0\"),\n    @(1, 5, \"64\u00f78=8, 0\"),\n\n    @(5, 1, \"71\u00f79=7, 8\"),\n    @(5, 2, \"32\u00f77=4, 4\"),\n    @(5, 3, \"85\u00f75=17, 0\"),\n    @(5, 4, \"86\u00f78=10, 6\"),\n    @(5, 5, \"37\u00f76=6, 1\"),\n\n    @(9, 1, \"44\u00f72=22, 0\"),\n    @(9, 2, \"34\u00f75=6, 4\"),\n    @(9, 3, \"27\u00f72=13, 1\"),\n    @(9, 4, \"17\u00f77=2, 3\"),\n    @(9, 5, \"34\u00f72=17, 0\"),\n\n    @(13, 1, \"93\u00f77=13, 2\"),\n    @(13, 2, \"86\u00f74=21, 2\"),\n    @(13, 3, \"82\u00f78=10, 2\"),\n    @(13, 4, \"13\u00f79=1, 4\"),\n    @(13, 5, \"52\u00f74=13, 0\"),\n\n    @(17, 1, \"86\u00f74=21, 2\"),\n    @(17, 2, \"47\u00f72=23, 1\"),\n    @(17, 3, \"70\u00f75=14, 0\"),\n    @(17, 4, \"16\u00f78=2, 0\"),\n    @(17, 5, \"66\u00f77=9, 3\")\n)\n\nforeach ($r in $replacements) {\n    $t.Cell($r[0], $r[1]).Range.Text = $r[2]\n}\n"}
